$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.185.03"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "3.935.04"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'611.42"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'171.63"
$ws.Range("D7").Value = "3.933.68"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'6.45"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'38.70"
$ws.Range("E13").Value = "  +5.29%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000258"
$ws.Range("E14").Value = "  +5.86%  "
$ws.Range("D15").Value = "4.598.50"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "3.951.63"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "70.204.98"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "'7.69"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "'18.61"
$ws.Range("E19").Value = "  +8.47%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").Value = "'497.66"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("D24").Value = "'0.0000166"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").Value = "'86.10"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "'10.23"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("D32").Value = "4.088.91"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "'7.91"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'32.36"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "3.901.10"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'6.18"
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "'3.30"
$ws.Range("E40").Value = "  +10.97%  "
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +8.49%  "
$ws.Range("D44").Value = "'439.07"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'48.38"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'8.70"
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0369"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000276"
$ws.Range("E49").Value = "  +22.30%  "
$ws.Range("D50").Value = "'40.79"
$ws.Range("E50").Value = "  +5.78%  "
$ws.Range("D51").Value = "'143.47"
$ws.Range("E51").Value = "  +0.52%  "
